$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": add a second (2014) column of source info in column K,
# plus three new footnote lines at the bottom explaining the CH4/N2O ratio
# methodology.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Re-word the existing 2010 citation block (B3 header text changes from
# "GHG Emissions" to "GHG Emissions (2010)").
$about.Range("B3").Value = "GHG Emissions (2010)"

# New column K mirrors column B's citation block, but for the 2014 data.
$about.Range("K3").Value = "GHG Emissions (2014)"
$about.Range("K4").Value = "Ministry of Environment, Forest and Climate Change"
$about.Range("K5").Value = "Second Biennial Update Report to the UNFCCC"
$about.Range("K6").Value = "https://unfccc.int/sites/default/files/resource/INDIA%20SECOND%20BUR%20High%20Res.pdf"
$about.Range("K7").Value = "Table 2.2: Greenhouse gas emissions, by sectors, for India in 2014"

# Copy the formatting from column B onto column K so it matches (bold+fill
# header, left-aligned source line, hyperlink style for the URL).
$about.Range("B3").Copy()
$about.Range("K3").PasteSpecial(-4122)
$about.Range("B4").Copy()
$about.Range("K4").PasteSpecial(-4122)
$about.Range("B6").Copy()
$about.Range("K6").PasteSpecial(-4122)

$about.Columns.Item(11).ColumnWidth = 57

# New footnote lines about the LULUCF-only CH4/N2O ratio methodology.
$about.Range("B16").Value = "For India, historical CH4/N2O emissions are available for the LULUCF sector only"
$about.Range("B17").Value = "for 2012 & 2014, in the Biennial reports. We use the ratios from the same to"
$about.Range("B18").Value = "estimate the average value to be applied to future years."

# ---------------------------------------------------------------------------
# Sheet "Data": fix the sign on the 2010 Net CO2 Emissions source value,
# and add a second table (rows 8-13) holding the 2014 Biennial Report data.
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

# Re-word the 2010 table header to disambiguate it from the new 2014 one.
$data.Range("A1").Value = "Table 2.2: GHG Emissions by sector (Gg) (First Biennial Report - 2010)"

# Fix sign error in the sourced CO2e figure for 2010.
$data.Range("F3").Value = -252531.78

# New 2014 table header (row 8), formatted like row 1.
$data.Range("A1").Copy()
$data.Range("A8").PasteSpecial(-4122)
$data.Range("B1").Copy()
$data.Range("B8:F8").PasteSpecial(-4122)
$data.Range("A8").Value = "Table 2.2: GHG Emissions by sector (Gg) (Second Biennial Report - 2014)"

# Column headers (row 9), same text as row 2.
$data.Range("B9").Value = "CO2 Emissions"
$data.Range("C9").Value = "CO2 Removals "
$data.Range("D9").Value = "CH4"
$data.Range("E9").Value = "N20"
$data.Range("F9").Value = "CO2e"

# 2014 LULUCF data row (row 10).
$data.Range("A10").Value = "LULUCF"
$data.Range("B10").Value = 17216.04
$data.Range("C10").Value = 319860.23
$data.Range("D10").Value = 48.19
$data.Range("E10").Value = 1.42
$data.Range("F10").Value = -301192.69

$data.Range("B10:C10").NumberFormat = "0.00"
$data.Range("F3").Copy()
$data.Range("F10").PasteSpecial(-4122)
$data.Range("F10").Value = -301192.69

# Net CO2 Emissions label + formula (rows 12-13), mirroring rows 5-6.
$data.Range("B12").Value = "Net CO2 Emissions"
$data.Range("B13").Formula = "=B10-C10"
$data.Range("B13").NumberFormat = "0.00"

$data.Columns.Item(6).ColumnWidth = 9.92

# ---------------------------------------------------------------------------
# Sheet "Calculations": add a parallel 2014 (Gg) column (C) next to the
# existing 2010 (Gg) column (B).
# ---------------------------------------------------------------------------
$calc = $wb.Worksheets.Item("Calculations")

$calc.Range("A1").Copy()
$calc.Range("C1").PasteSpecial(-4122)

$calc.Range("B3").Value = "2010 (Gg)"
$calc.Range("C3").Value = "2014 (Gg)"

$calc.Range("C4").Formula = "=Data!B13"
$calc.Range("C4").NumberFormat = "0.00"

$calc.Range("C5").Formula = "=Data!D10"
$calc.Range("C6").Formula = "=Data!E10"

$calc.Range("C8").Formula = "=C5/C4"
$calc.Range("B8").Copy()
$calc.Range("C8").PasteSpecial(-4122)
$calc.Range("C8").Formula = "=C5/C4"

$calc.Range("C9").Formula = "=C6/C4"
$calc.Range("B9").Copy()
$calc.Range("C9").PasteSpecial(-4122)
$calc.Range("C9").Formula = "=C6/C4"

$calc.Columns.Item(3).ColumnWidth = 9.6

# ---------------------------------------------------------------------------
# Sheet "RPEpUACE": the CH4/CO2 and N2O/CO2 ratios now average the 2010 and
# 2014 estimates (negated, since the source ratios come out negative).
# ---------------------------------------------------------------------------
$rpe = $wb.Worksheets.Item("RPEpUACE")
$rpe.Range("B11").Formula = "=-AVERAGE(Calculations!B8, Calculations!C8)"
$rpe.Range("B12").Formula = "=-AVERAGE(Calculations!B9,Calculations!C9)"

# Update selections to match the final authored state, and make RPEpUACE the
# active (visible) tab.
$null = $about.Range("B13").Select()
$null = $data.Range("B12").Select()
$null = $calc.Range("B8").Select()
$null = $rpe.Range("B15").Select()
$rpe.Activate()
